$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Förändrad" column (C) holds a date serial that was bumped by one day
# (2026-02-21 -> 2026-02-22, i.e. 46074 -> 46075) for every data row (2..340).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 340 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value -ne $null) {
        $cell.Value = 46075
    }
}
